# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (commit: "Updated cryptos list on Sun Dec 31 11:34:55 UTC
# 2023 with GitHub Actions"). All target cells in this sheet hold their
# data as literal text (OOXML inlineStr / shared string), even when the
# text looks like a number (e.g. "318.53" or "42.766.05" - note some
# "prices" even contain two dots and are never valid numbers). Assigning
# a plain numeric-looking string straight to Range.Value lets Excel's
# COM layer auto-convert it into a real number, which would corrupt
# both the stored value (e.g. "22.70" -> 22.7) and the cell's type/style.
# To avoid that we force each cell to Text format before writing, then
# restore the original "Normal" style (so no stray number formats are
# left behind) while keeping the freshly written text intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '42.766.05' },
    @{ Cell = 'E2'; Value = '  +1.69%  ' },
    @{ Cell = 'D3'; Value = '2.311.99' },
    @{ Cell = 'E3'; Value = '  +1.09%  ' },
    @{ Cell = 'E4'; Value = '  -0.03%  ' },
    @{ Cell = 'D5'; Value = '318.53' },
    @{ Cell = 'E5'; Value = '  +0.59%  ' },
    @{ Cell = 'D6'; Value = '104.56' },
    @{ Cell = 'E6'; Value = '  +1.84%  ' },
    @{ Cell = 'D7'; Value = '0.631' },
    @{ Cell = 'E7'; Value = '  +0.78%  ' },
    @{ Cell = 'E8'; Value = '  +0.05%  ' },
    @{ Cell = 'D9'; Value = '0.608' },
    @{ Cell = 'E9'; Value = '  +0.94%  ' },
    @{ Cell = 'D10'; Value = '40.14' },
    @{ Cell = 'E10'; Value = '  +2.54%  ' },
    @{ Cell = 'D11'; Value = '0.0909' },
    @{ Cell = 'E11'; Value = '  +0.53%  ' },
    @{ Cell = 'D12'; Value = '8.55' },
    @{ Cell = 'E12'; Value = '  +3.75%  ' },
    @{ Cell = 'E13'; Value = '  +0.75%  ' },
    @{ Cell = 'D14'; Value = '0.978' },
    @{ Cell = 'E14'; Value = '  +1.82%  ' },
    @{ Cell = 'D15'; Value = '15.41' },
    @{ Cell = 'E15'; Value = '  +1.27%  ' },
    @{ Cell = 'D16'; Value = '2.663.84' },
    @{ Cell = 'E16'; Value = '  +1.11%  ' },
    @{ Cell = 'D17'; Value = '2.334.41' },
    @{ Cell = 'E17'; Value = '  +2.29%  ' },
    @{ Cell = 'D18'; Value = '42.731.25' },
    @{ Cell = 'E18'; Value = '  +1.74%  ' },
    @{ Cell = 'D19'; Value = '7.52' },
    @{ Cell = 'E19'; Value = '  +0.78%  ' },
    @{ Cell = 'E20'; Value = '  +1.08%  ' },
    @{ Cell = 'D21'; Value = '13.31' },
    @{ Cell = 'E21'; Value = '  +32.81%  ' },
    @{ Cell = 'D22'; Value = '73.83' },
    @{ Cell = 'E22'; Value = '  +0.73%  ' },
    @{ Cell = 'E23'; Value = '  -2.18%  ' },
    @{ Cell = 'D24'; Value = '268.82' },
    @{ Cell = 'E24'; Value = '  -4.88%  ' },
    @{ Cell = 'E25'; Value = '  +0.40%  ' },
    @{ Cell = 'E26'; Value = '  -0.52%  ' },
    @{ Cell = 'E27'; Value = '  +1.22%  ' },
    @{ Cell = 'D29'; Value = '22.70' },
    @{ Cell = 'E29'; Value = '  -0.84%  ' },
    @{ Cell = 'D30'; Value = '37.94' },
    @{ Cell = 'E30'; Value = '  +6.30%  ' },
    @{ Cell = 'B31'; Value = 'Monero' },
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D31'; Value = '165.82' },
    @{ Cell = 'E31'; Value = '  +1.53%  ' },
    @{ Cell = 'B32'; Value = 'Filecoin' },
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Cell = 'D32'; Value = '6.22' },
    @{ Cell = 'E32'; Value = '  +6.66%  ' },
    @{ Cell = 'D33'; Value = '0.0892' },
    @{ Cell = 'E33'; Value = '  +2.16%  ' },
    @{ Cell = 'E34'; Value = '  -1.32%  ' },
    @{ Cell = 'D35'; Value = '2.58' },
    @{ Cell = 'E35'; Value = '  -9.40%  ' },
    @{ Cell = 'E36'; Value = '  +0.87%  ' },
    @{ Cell = 'E37'; Value = '  +2.37%  ' },
    @{ Cell = 'E38'; Value = '  +1.54%  ' },
    @{ Cell = 'E39'; Value = '  +1.59%  ' },
    @{ Cell = 'E40'; Value = '  -4.46%  ' },
    @{ Cell = 'E41'; Value = '  +10.49%  ' },
    @{ Cell = 'D42'; Value = '98.86' },
    @{ Cell = 'E42'; Value = '  -1.29%  ' },
    @{ Cell = 'D43'; Value = '70.20' },
    @{ Cell = 'E43'; Value = '  +1.18%  ' },
    @{ Cell = 'E44'; Value = '  +1.05%  ' },
    @{ Cell = 'E45'; Value = '  -0.13%  ' },
    @{ Cell = 'D46'; Value = '12.39' },
    @{ Cell = 'E46'; Value = '  +4.52%  ' },
    @{ Cell = 'D47'; Value = '82.28' },
    @{ Cell = 'E47'; Value = '  +7.59%  ' },
    @{ Cell = 'D48'; Value = '115.11' },
    @{ Cell = 'E48'; Value = '  +0.81%  ' },
    @{ Cell = 'B49'; Value = 'FraxShare' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Cell = 'D49'; Value = '8.91' },
    @{ Cell = 'E49'; Value = '  -0.34%  ' },
    @{ Cell = 'B50'; Value = 'THORChain' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' },
    @{ Cell = 'D50'; Value = '5.30' },
    @{ Cell = 'E50'; Value = '  +0.64%  ' },
    @{ Cell = 'D51'; Value = '1.620.57' },
    @{ Cell = 'E51'; Value = '  +5.38%  ' }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $cell.NumberFormat = '@'
    $cell.Value = $update.Value
    $cell.Style = 'Normal'
}
